$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 4 (shifts "Project Descriptions" down to row 5)
$ws.Rows.Item(4).Insert()

# Update A2 to combined label
$ws.Range("A2").Value = "Client/Coder Sign-in"

# A3 becomes "Client Registration" (was "Coder Sign-in" before insert, now overwritten)
$ws.Range("A3").Value = "Client Registration"

# New row 4 gets "Coder Registration"
$ws.Range("A4").Value = "Coder Registration"

# Row5 already holds "Project Descriptions" (shifted down automatically by Insert)

$ws.Range("A4").Select()
